$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Pred_pCO2 formula (AD7) to use the re-indexed model fit coefficients.
$ws.Range("AD7").Formula = "=10^(1.553889386876 + 1.151209200966*Z7 + 0.028619316030*T7 + -0.055064417137*Z7^2 + -0.018467383046*Z7*T7 + -0.000532223052*T7^2 + 0.008443439468*Z7^2*T7 + 0.000734604699*Z7*T7^2 + -0.000336681023*Z7^2*T7^2)"

# Add new row that documents the coefficient equation used above.
$ws.Range("T9").Value = "log10(pCO2) = "
$ws.Range("U9").Value = "(1.553889386876 +  1.151209200966*Z7 +  0.028619316030*T7 +  -0.055064417137*Z7^2 +  -0.018467383046*Z7*T7 +  -0.000532223052*T7^2 +  0.008443439468*Z7^2*T7 +  0.000734604699*Z7*T7^2 +  -0.000336681023*Z7^2*T7^2)"

# Widen columns T, Z, and AD to fit the new content (~15.71 chars).
$ws.Columns.Item(20).ColumnWidth = 14.83
$ws.Columns.Item(26).ColumnWidth = 14.83
$ws.Columns.Item(30).ColumnWidth = 14.83
